$d = $word.ActiveDocument

$para1 = "BOG BOB BOGGY BOG BOG BOG BOB BOGGY BOG BOG"
$para2 = "BOG BOB BOGGY BOG BOG"

$r = $d.Paragraphs(1).Range
$r.InsertBefore($para2 + "`r")
$r.InsertBefore($para2 + "`r")
$r.InsertBefore($para2 + "`r")
$r.InsertBefore($para2 + "`r")
$r.InsertBefore($para1 + "`r")

$count = $d.Paragraphs.Count
Write-Output "paragraph count: $count"

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $p.Range.Font.Bold = 1
    $p.Range.Font.Size = 72
}
